# Add toolbar options and update README.md
#
# Puts the text "dsfsdfsd" (bold) into cell A1 of the active sheet, and
# switches the page setup to portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "dsfsdfsd"
$ws.Range("A1").Font.Bold = $true

# xlPortrait = 1
$ws.PageSetup.Orientation = 1
